# Insert 3 new weekly-record rows at the top of the data block (rows 260-262),
# pushing the existing rows 260-276 down to 263-279, then populate the new
# rows with this week's values (date 2021-11-16 / serial 44516).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows 260..276 down by three rows.
$ws.Rows.Item(260).Resize(3).Insert()

# --- Row 260 ---
$ws.Cells.Item(260, 1).Value = 10
$ws.Cells.Item(260, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(260, 3).Value = "La Araucanía"
$ws.Cells.Item(260, 4).Value = 44516
$ws.Cells.Item(260, 5).Value = 9
$ws.Cells.Item(260, 6).Value = 100112008
$ws.Cells.Item(260, 7).Value = "Coliflor"
$ws.Cells.Item(260, 8).Value = "Sin especificar"
$ws.Cells.Item(260, 9).Value = "Primera"
$ws.Cells.Item(260, 10).Value = 750
$ws.Cells.Item(260, 11).Value = 900
$ws.Cells.Item(260, 12).Value = 900
$ws.Cells.Item(260, 13).Value = 900
$ws.Cells.Item(260, 14).Value = "`$/unidad"
$ws.Cells.Item(260, 15).Value = "Región Metropolitana"
$ws.Cells.Item(260, 16).Value = 900
$ws.Cells.Item(260, 17).Value = 1
$ws.Cells.Item(260, 18).Value = "Hortaliza"

# --- Row 261 ---
$ws.Cells.Item(261, 1).Value = 10
$ws.Cells.Item(261, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(261, 3).Value = "La Araucanía"
$ws.Cells.Item(261, 4).Value = 44516
$ws.Cells.Item(261, 5).Value = 9
$ws.Cells.Item(261, 6).Value = 100112008
$ws.Cells.Item(261, 7).Value = "Coliflor"
$ws.Cells.Item(261, 8).Value = "Sin especificar"
$ws.Cells.Item(261, 9).Value = "Primera"
$ws.Cells.Item(261, 10).Value = 1500
$ws.Cells.Item(261, 11).Value = 800
$ws.Cells.Item(261, 12).Value = 900
$ws.Cells.Item(261, 13).Value = 850
$ws.Cells.Item(261, 14).Value = "`$/unidad"
$ws.Cells.Item(261, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(261, 16).Value = 850
$ws.Cells.Item(261, 17).Value = 1
$ws.Cells.Item(261, 18).Value = "Hortaliza"

# --- Row 262 ---
$ws.Cells.Item(262, 1).Value = 10
$ws.Cells.Item(262, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(262, 3).Value = "La Araucanía"
$ws.Cells.Item(262, 4).Value = 44516
$ws.Cells.Item(262, 5).Value = 9
$ws.Cells.Item(262, 6).Value = 100112008
$ws.Cells.Item(262, 7).Value = "Coliflor"
$ws.Cells.Item(262, 8).Value = "Sin especificar"
$ws.Cells.Item(262, 9).Value = "Primera"
$ws.Cells.Item(262, 10).Value = 2800
$ws.Cells.Item(262, 11).Value = 800
$ws.Cells.Item(262, 12).Value = 900
$ws.Cells.Item(262, 13).Value = 845
$ws.Cells.Item(262, 14).Value = "`$/unidad"
$ws.Cells.Item(262, 15).Value = "Región del Maule"
$ws.Cells.Item(262, 16).Value = 845
$ws.Cells.Item(262, 17).Value = 1
$ws.Cells.Item(262, 18).Value = "Hortaliza"
